# Auto update stock data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that receive a new date (all moving from 2026/01/15 to 2026/01/16)
# Column A is stored as plain text (not a real date), so force the Text
# number format first so Excel does not auto-convert the string into a
# date serial value, then restore the original "Normal" cell style so no
# formatting is left behind on the cell.
$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = "2026/01/16"
    $cell.Style = "Normal"
}

# Rows with an updated EBITDA value (column B) - also stored as text
$bvals = @{ 2 = "7.97"; 8 = "8.95"; 14 = "3.25"; 20 = "14.04"; 26 = "12.14"; 32 = "29.75"; 44 = "17.00"; 50 = "12.76"; 56 = "32.80"; 62 = "12.50"; 68 = "13.80"; 74 = "19.86" }
foreach ($r in $bvals.Keys) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $bvals[$r]
    $cell.Style = "Normal"
}
